# Edit excel data: add new user info (update username/email on row 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B3 with the new email address (replacing the previous one)
$ws.Range("B3").Value = "ogulcan.a81@hotmail.com"

# Update the active selection to B3, matching the saved view state
$ws.Range("B3").Select()
